# SlicerLibraries.xlsx update — "updates to Eigen 3.4.0 3147391d"
#
# The sheet tracks third-party library versions. This commit:
#   - fills in the "boost nowide" row: Version_new = 11.3.0, status = ok
#   - fills in the "eigen" row: Version old = 3.3.7 323c052e1731, Version_new = 3.4.0
#   - fills in the "exif" row: Version_new = 2002 (matches Version old)
#   - normalizes the "status" column formatting for the whole lower block
#     (rows 30-57) to match the formatting already used higher up the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- content updates -------------------------------------------------

# " boost nowide" row (row 35): new version + status
$ws.Range("C35").Value = "11.3.0"
$ws.Range("D35").Value = "ok"

# " eigen" row (row 38): old version (with commit) -> new version
$ws.Range("B38").Value = "3.3.7 323c052e1731"
$ws.Range("C38").Value = "3.4.0"

# " exif" row (row 39): Version_new mirrors Version old
$ws.Range("C39").Value = 2002

# --- formatting normalization -----------------------------------------
# D30:D57 currently use a slightly different (duplicate) "status" cell
# style than the rest of column D. Bring them in line by copying the
# format already used in D2 over the whole block.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D30:D57").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$excel.CutCopyMode = $false

# --- restore selection -------------------------------------------------
$ws.Range("D38").Select() | Out-Null
